$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# NOTE: cells are edited in row 4 -> row 3 -> row 2 order so that the
# workbook's shared-string table is rebuilt in the same order the source
# workbook uses (Files query, then Samples query, then Case ID query).

# --- Row 4: Files query (column B) - replace trailing "order by f.file_name" with
#     "order By f.file_name ASC LIMIT 100" (also fixes capitalisation of "By") ---
$b4 = $ws.Range("B4").Value2
$b4 = $b4 -replace "    order by f\.file_name$", "     order By f.file_name ASC LIMIT 100"
$ws.Range("B4").Value2 = $b4

# --- Row 3: Sample ID query (column B) gets an ORDER BY / LIMIT clause appended ---
$b3 = $ws.Range("B3").Value2
$ws.Range("B3").Value2 = $b3 + "`n order By samp.sample_id ASC LIMIT 100"

# --- Row 2: Case ID query (column B) gets an ORDER BY / LIMIT clause appended ---
$b2 = $ws.Range("B2").Value2
$ws.Range("B2").Value2 = $b2 + "`n order By ss.study_subject_id ASC LIMIT 100"

# --- Row heights grow by one wrapped line for rows 2 & 3 (row 4 already at the cap) ---
$ws.Rows(2).RowHeight = 331.2
$ws.Rows(3).RowHeight = 360

# --- Selection moves from B4 to B2 ---
$ws.Range("B2").Select()
